# Generate Report for Handback
# Update the timestamp strings recorded during handback report generation.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date for the first file.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-04 15:09:12"

# "zh-cn" sheet: Correspond Handoff/Handback Datetime for the first file.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-04 15:09:08"
$wsZhCn.Range("K2").Value = "2016-09-04 15:09:26"

# "de-de" sheet: Correspond Handoff/Handback Datetime for the first file.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-04 15:09:12"
$wsDeDe.Range("K2").Value = "2016-09-04 15:09:33"
